# Insert a new "chemical_recycling_pyrolysis" parameter row directly
# below the existing "chemical_recycling_gasification" row (row 9),
# pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10 (shifts rows 10..24 down to 11..25)
$ws.Rows.Item(10).Insert()

# Populate the new row with the new parameter and its value
$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
